# Applies the "Added a few more slots" edit to the Blood Lust review document:
#  1. Inserts a new "Meta description" paragraph right after the title heading.
#  2. Removes the duplicated bold title paragraph near the end of the document.
#  3. Replaces the remaining (italic) paragraph's text with a DALL-E image prompt.

$d = $word.ActiveDocument

function New-WordPackageXml([string]$bodyInnerXml) {
    return '<?xml version="1.0" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body>' +
    $bodyInnerXml +
    '<w:sectPr><w:pgSz w:w="12240" w:h="15840"/></w:sectPr>' +
    '</w:body></w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'
}

# --- Step 1: work at the end of the document first, so paragraph indices  ---
# --- for the earlier part of the document are not disturbed.             ---

# The last two paragraphs are, in order:
#   50: "Limited betting options"            (unrelated bullet item)
#   51: bold  "Play Blood Lust Free Slot | Review & Guide 2021"
#   52: italic "Read our unbiased review ..."
$lastIndex = $d.Paragraphs.Count
$boldTitlePara = $d.Paragraphs.Item($lastIndex - 1)
Write-Host "Removing paragraph: [" $boldTitlePara.Range.Text "]"
$boldTitlePara.Range.Delete()

# Replace the now-last paragraph's text (still italic) with the DALL-E prompt,
# keeping the leading empty run / italic run structure intact.
$dallePrompt = 'DALLE, please create a feature image fitting the game "Blood Lust". The image should be in a cartoon style and feature a happy Maya warrior with glasses. Ensure that the image is suitable for online slot games and is visually appealing to the target audience. You may use elements from the game, such as the thematic icons and dark background colors, to enhance the image and capture the essence of Blood Lust.'

$italicPara = $d.Paragraphs.Item($d.Paragraphs.Count)
Write-Host "Replacing paragraph: [" $italicPara.Range.Text "]"
$italicParaXml = New-WordPackageXml ('<w:p><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>{0}</w:t></w:r></w:p>' -f $dallePrompt)
$italicPara.Range.InsertXML($italicParaXml)

# --- Step 2: insert the new "Meta description" paragraph after the title. ---

$metaDescriptionXml = New-WordPackageXml '<w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t xml:space="preserve">: Read our unbiased review of Blood Lust video slot, learn bonus features, tips to play and win. Play Blood Lust slot online free without download.</w:t></w:r></w:p>'

$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.InsertParagraphAfter()
$newMetaPara = $d.Paragraphs.Item(2)
$newMetaPara.Range.InsertXML($metaDescriptionXml)

Write-Host "Final paragraph count: " $d.Paragraphs.Count
